# "Refined metadata to be additional tab"
# Adds a new "metadata" worksheet (after "data") summarising the panel
# query, and refreshes the "time_taken" timestamps on the "data" sheet
# to the time of the (re-)fetch.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Refresh the per-row query timestamps on the "data" sheet (column F).
$dataSheet.Range("F2").Value = "2021-10-05 14:21:29.807991"
$dataSheet.Range("F3").Value = "2021-10-05 14:21:29.807999"
$dataSheet.Range("F4").Value = "2021-10-05 14:21:29.808002"
$dataSheet.Range("F5").Value = "2021-10-05 14:21:29.808004"
$dataSheet.Range("F6").Value = "2021-10-05 14:21:29.808007"
$dataSheet.Range("F7").Value = "2021-10-05 14:21:29.808009"
$dataSheet.Range("F8").Value = "2021-10-05 14:21:29.808012"
$dataSheet.Range("F9").Value = "2021-10-05 14:21:29.808014"
$dataSheet.Range("F10").Value = "2021-10-05 14:21:29.808017"

# Add the new "metadata" sheet right after "data".
$meta = $wb.Worksheets.Add($null, $dataSheet)
$meta.Name = "metadata"

# Header row (B1:G1) — same bold/bordered header style as "data"'s header.
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$dataSheet.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

$dataSheet.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data row (A2:G2) describing the panel this workbook was built from.
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Membranoproliferative glomerulonephritis"
$meta.Range("C2").Value = 83

# data_version ("2.21") must stay text, not be coerced to a number.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "2.21"
$meta.Range("D2").ClearFormats()

$meta.Range("E2").Value = "2021-07-13T10:34:25.389534Z"
$meta.Range("F2").Value = "2021-10-05 14:21:29.805063"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/83/?format=json"

# Leave "data" as the active/selected sheet, as before the edit.
$dataSheet.Activate() | Out-Null
$dataSheet.Range("A1").Select() | Out-Null

Write-Host "done"
